$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = -0.0192;  C = -0.0192;  D = -0.03030000068247318; E = 0.2311999946832657; F = -0.07940000295639038 }
    3 = @{ B = -0.0562;  C = -0.0562;  D = -0.0794;               E = 0.2745999991893768; F = -0.07940000295639038 }
    4 = @{ B = -0.0658;  C = -0.0658;  D = -0.0794;               E = 0.2436999976634979; F = -0.07940000295639038 }
    5 = @{ B = 0.0953;   C = 0.0953;   D = 0.1224;                E = 0.354200005531311;  F = -0.07940000295639038 }
    6 = @{ B = 0.0742;   C = 0.0742;   D = 0.0531;                E = 0.3265999853610992; F = -0.07940000295639038 }
    7 = @{ B = 0.023;    C = 0.023;    D = 0.01250000018626451;   E = 0.221000000834465;  F = -0.06440000236034393 }
    8 = @{ B = 0.027;    C = 0.027;    D = 0.0179;                E = 0.1738000065088272; F = -0.06440000236034393 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
